$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H114").Value = 70722
$ws.Range("J114").Value = 70722
$ws.Range("L114").Value = 70722
$ws.Range("N114").Value = -79400

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9973.096
$ws.Range("I32").Value = 8847.192999999999
$ws.Range("K32").Value = 8847.192999999999
$ws.Range("M32").Value = -8560.192999999999

$ws.Range("H110").Value = 1796.5769
$ws.Range("I110").Value = 1703.6666
$ws.Range("J110").Value = 2911.5
$ws.Range("K110").Value = 1703.6666
$ws.Range("L110").Value = 2911.5
$ws.Range("M110").Value = 341.3334
$ws.Range("N110").Value = -7001.5

$ws.Range("H135").Value = 36109.332
$ws.Range("J135").Value = 36109.332
$ws.Range("L135").Value = 36109.332
$ws.Range("N135").Value = -46249.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2672.818
$ws.Range("I20").Value = 2325
$ws.Range("J20").Value = 2871.5715
$ws.Range("K20").Value = 2325
$ws.Range("L20").Value = 2871.5715
$ws.Range("M20").Value = -2078
$ws.Range("N20").Value = -3365.5715

$ws.Range("H86").Value = 1808.25
$ws.Range("I86").Value = 1868.1666
$ws.Range("J86").Value = 1628.5
$ws.Range("K86").Value = 1868.1666
$ws.Range("L86").Value = 1628.5
$ws.Range("M86").Value = -745.1666
$ws.Range("N86").Value = -3874.5

$ws.Range("H89").Value = 1808.25
$ws.Range("I89").Value = 1868.1666
$ws.Range("J89").Value = 1628.5
$ws.Range("K89").Value = 9340.833000000001
$ws.Range("L89").Value = 8142.5
$ws.Range("M89").Value = -3724.833000000001
$ws.Range("N89").Value = -19374.5

$ws.Range("H94").Value = 1920.0588
$ws.Range("I94").Value = 1828.7916
$ws.Range("J94").Value = 2139.1
$ws.Range("K94").Value = 1828.7916
$ws.Range("L94").Value = 2139.1
$ws.Range("M94").Value = -1377.7916
$ws.Range("N94").Value = -3041.1

$ws.Range("H134").Value = 2045.3214
$ws.Range("I134").Value = 1766.8636
$ws.Range("J134").Value = 3066.3333
$ws.Range("K134").Value = 5300.5908
$ws.Range("L134").Value = 9198.999899999999
$ws.Range("M134").Value = -2765.5908
$ws.Range("N134").Value = -14268.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4698799
$ws.Range("I31").Value = 1556.7838
$ws.Range("J31").Value = 9810504
$ws.Range("K31").Value = 1556.7838
$ws.Range("L31").Value = 9810504
$ws.Range("M31").Value = -1261.7838
$ws.Range("N31").Value = -9811094

$ws.Range("H34").Value = 4698799
$ws.Range("I34").Value = 1556.7838
$ws.Range("J34").Value = 9810504
$ws.Range("K34").Value = 1556.7838
$ws.Range("L34").Value = 9810504
$ws.Range("M34").Value = -1354.7838
$ws.Range("N34").Value = -9810908

$ws.Range("H99").Value = 3788.5833
$ws.Range("I99").Value = 4007.6667
$ws.Range("J99").Value = 3715.5557
$ws.Range("K99").Value = 4007.6667
$ws.Range("L99").Value = 3715.5557
$ws.Range("M99").Value = -2509.6667
$ws.Range("N99").Value = -6711.5557

$ws.Range("H126").Value = 3788.5833
$ws.Range("I126").Value = 4007.6667
$ws.Range("J126").Value = 3715.5557
$ws.Range("K126").Value = 12023.0001
$ws.Range("L126").Value = 11146.6671
$ws.Range("M126").Value = -9553.000100000001
$ws.Range("N126").Value = -16086.6671

$ws.Range("H134").Value = 501394.4
$ws.Range("I134").Value = 543085.25
$ws.Range("J134").Value = 284602
$ws.Range("K134").Value = 1629255.75
$ws.Range("L134").Value = 853806
$ws.Range("M134").Value = -1626720.75
$ws.Range("N134").Value = -858876

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2802.754
$ws.Range("I131").Value = 11534.333
$ws.Range("J131").Value = 1399.4642
$ws.Range("K131").Value = 34602.999
$ws.Range("L131").Value = 4198.392599999999
$ws.Range("M131").Value = -29562.999
$ws.Range("N131").Value = -14278.3926

$ws.Range("H134").Value = 64377840
$ws.Range("I134").Value = 79233224
$ws.Range("J134").Value = 4500
$ws.Range("K134").Value = 237699672
$ws.Range("L134").Value = 13500
$ws.Range("M134").Value = -237694602
$ws.Range("N134").Value = -23640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5981.6665
$ws.Range("I80").Value = 5127.8
$ws.Range("J80").Value = 8650
$ws.Range("K80").Value = 5127.8
$ws.Range("L80").Value = 8650
$ws.Range("M80").Value = -4129.8
$ws.Range("N80").Value = -10646

$ws.Range("H83").Value = 5981.6665
$ws.Range("I83").Value = 5127.8
$ws.Range("J83").Value = 8650
$ws.Range("K83").Value = 25639
$ws.Range("L83").Value = 43250
$ws.Range("M83").Value = -20647
$ws.Range("N83").Value = -53234

$ws.Range("H122").Value = 1590.9445
$ws.Range("I122").Value = 1666.9286
$ws.Range("J122").Value = 1325
$ws.Range("K122").Value = 5000.7858
$ws.Range("L122").Value = 3975
$ws.Range("M122").Value = -2550.7858
$ws.Range("N122").Value = -8875

$ws.Range("H126").Value = 9224.933999999999
$ws.Range("I126").Value = 17071.428
$ws.Range("J126").Value = 2359.25
$ws.Range("K126").Value = 51214.284
$ws.Range("L126").Value = 7077.75
$ws.Range("M126").Value = -48744.284
$ws.Range("N126").Value = -12017.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H16").Value = 2803
$ws.Range("I16").Value = 1830.5454
$ws.Range("J16").Value = 13500
$ws.Range("K16").Value = 1830.5454
$ws.Range("L16").Value = 13500
$ws.Range("M16").Value = -1660.5454
$ws.Range("N16").Value = -13840

$ws.Range("H46").Value = 9458.444
$ws.Range("I46").Value = 1090
$ws.Range("J46").Value = 12677.077
$ws.Range("K46").Value = 1090
$ws.Range("L46").Value = 12677.077
$ws.Range("M46").Value = -902
$ws.Range("N46").Value = -13053.077

$ws.Range("H93").Value = 894.9091
$ws.Range("I93").Value = 549.1429000000001
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 549.1429000000001
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = 698.8570999999999
$ws.Range("N93").Value = -3996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 45000
$ws.Range("J116").Value = 45000
$ws.Range("L116").Value = 45000
